{"js": "// Insert \" after a random amount of time\" right before \" (LoraWan-like)\"\n// in the ack-confirmation paragraph, i.e. turn:\n//   \"...it sends again the message (LoraWan-like). ...\"\n// into:\n//   \"...it sends again the message after a random amount of time (LoraWan-like). ...\"\n\nconst results = context.document.body.search(\"it sends again the message\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  target.insertText(\" after a random amount of time\", Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "# Insert \" after a random amount of time\" right before \" (LoraWan-like)\"\n# in the ack-confirmation paragraph, i.e. turn:\n#   \"...it sends again the message (LoraWan-like). ...\"\n# into:\n#   \"...it sends again the message after a random amount of time (LoraWan-like). ...\"\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"it sends again the message\")\n\nif ($found) {\n    # Collapse the found range to its end point (right after \"...the message\")\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.InsertAfter(\" after a random amount of time\")\n}\n"}
